# Scheduled runner update: refresh cached market-board derived values
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 5673.6665
$ws.Range("I20").Value = 5673.6665
$ws.Range("K20").Value = 5673.6665
$ws.Range("M20").Value = -5443.6665
$ws.Range("H33").Value = 250.14285
$ws.Range("I33").Value = 271.25
$ws.Range("J33").Value = 123.5
$ws.Range("K33").Value = 271.25
$ws.Range("L33").Value = 123.5
$ws.Range("M33").Value = -42.25
$ws.Range("N33").Value = -581.5
$ws.Range("H35").Value = 5673.6665
$ws.Range("I35").Value = 5673.6665
$ws.Range("K35").Value = 5673.6665
$ws.Range("M35").Value = -5294.6665
$ws.Range("H74").Value = 3900.2727
$ws.Range("I74").Value = 3145.75
$ws.Range("J74").Value = 4331.4287
$ws.Range("K74").Value = 3145.75
$ws.Range("L74").Value = 4331.4287
$ws.Range("M74").Value = -2209.75
$ws.Range("N74").Value = -6203.4287
$ws.Range("H77").Value = 3900.2727
$ws.Range("I77").Value = 3145.75
$ws.Range("J77").Value = 4331.4287
$ws.Range("K77").Value = 15728.75
$ws.Range("L77").Value = 21657.1435
$ws.Range("M77").Value = -11048.75
$ws.Range("N77").Value = -31017.1435
$ws.Range("H98").Value = 69445440
$ws.Range("I98").Value = 89286790
$ws.Range("J98").Value = 701
$ws.Range("K98").Value = 89286790
$ws.Range("L98").Value = 701
$ws.Range("M98").Value = -89285292
$ws.Range("N98").Value = -3697
$ws.Range("H113").Value = 2661.2144
$ws.Range("I113").Value = 2664.8572
$ws.Range("J113").Value = 2657.5715
$ws.Range("K113").Value = 2664.8572
$ws.Range("L113").Value = 2657.5715
$ws.Range("M113").Value = 589.1428000000001
$ws.Range("N113").Value = -9165.5715
$ws.Range("H122").Value = 69445440
$ws.Range("I122").Value = 89286790
$ws.Range("J122").Value = 701
$ws.Range("K122").Value = 267860370
$ws.Range("L122").Value = 2103
$ws.Range("M122").Value = -267857920
$ws.Range("N122").Value = -7003
$ws.Range("H137").Value = 1085.68
$ws.Range("I137").Value = 866.1053000000001
$ws.Range("J137").Value = 1781
$ws.Range("K137").Value = 2598.3159
$ws.Range("L137").Value = 5343
$ws.Range("M137").Value = -48.31590000000006
$ws.Range("N137").Value = -10443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3705247
$ws.Range("I61").Value = 4631012.5
$ws.Range("J61").Value = 2186
$ws.Range("K61").Value = 4631012.5
$ws.Range("L61").Value = 2186
$ws.Range("M61").Value = -4630800.5
$ws.Range("N61").Value = -2610
$ws.Range("H63").Value = 2791.0256
$ws.Range("I63").Value = 2860.889
$ws.Range("J63").Value = 2633.8333
$ws.Range("K63").Value = 2860.889
$ws.Range("L63").Value = 2633.8333
$ws.Range("M63").Value = -2174.889
$ws.Range("N63").Value = -4005.8333
$ws.Range("H66").Value = 2791.0256
$ws.Range("I66").Value = 2860.889
$ws.Range("J66").Value = 2633.8333
$ws.Range("K66").Value = 14304.445
$ws.Range("L66").Value = 13169.1665
$ws.Range("M66").Value = -10872.445
$ws.Range("N66").Value = -20033.1665
$ws.Range("H74").Value = 1095.6383
$ws.Range("I74").Value = 1094.875
$ws.Range("K74").Value = 1094.875
$ws.Range("M74").Value = -220.875
$ws.Range("H77").Value = 1095.6383
$ws.Range("I77").Value = 1094.875
$ws.Range("K77").Value = 5474.375
$ws.Range("M77").Value = -1106.375
$ws.Range("H132").Value = 1154774.4
$ws.Range("I132").Value = 1138.3513
$ws.Range("J132").Value = 4203669.5
$ws.Range("K132").Value = 3415.0539
$ws.Range("L132").Value = 12611008.5
$ws.Range("M132").Value = -885.0538999999999
$ws.Range("N132").Value = -12616068.5
$ws.Range("H136").Value = 3705247
$ws.Range("I136").Value = 4631012.5
$ws.Range("J136").Value = 2186
$ws.Range("K136").Value = 13893037.5
$ws.Range("L136").Value = 6558
$ws.Range("M136").Value = -13890487.5
$ws.Range("N136").Value = -11658

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1636488.4
$ws.Range("I134").Value = 1007.22644
$ws.Range("J134").Value = 7415188.5
$ws.Range("K134").Value = 3021.67932
$ws.Range("L134").Value = 22245565.5
$ws.Range("M134").Value = -486.6793200000002
$ws.Range("N134").Value = -22250635.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1278
$ws.Range("H31").Value = 1049520.5
$ws.Range("I31").Value = 1588165.4
$ws.Range("J31").Value = 2155.5
$ws.Range("K31").Value = 1588165.4
$ws.Range("L31").Value = 2155.5
$ws.Range("M31").Value = -1587870.4
$ws.Range("N31").Value = -2745.5
$ws.Range("H32").Value = 4759.8
$ws.Range("I32").Value = 3449.75
$ws.Range("K32").Value = 3449.75
$ws.Range("M32").Value = -3133.75
$ws.Range("H34").Value = 1049520.5
$ws.Range("I34").Value = 1588165.4
$ws.Range("J34").Value = 2155.5
$ws.Range("K34").Value = 1588165.4
$ws.Range("L34").Value = 2155.5
$ws.Range("M34").Value = -1587963.4
$ws.Range("N34").Value = -2559.5
$ws.Range("H58").Value = 28572172
$ws.Range("I58").Value = 52632424
$ws.Range("J58").Value = 625.5
$ws.Range("K58").Value = 52632424
$ws.Range("L58").Value = 625.5
$ws.Range("M58").Value = -52632221
$ws.Range("N58").Value = -1031.5
$ws.Range("H132").Value = 9525362
$ws.Range("I132").Value = 1189.421
$ws.Range("J132").Value = 20835316
$ws.Range("K132").Value = 3568.263
$ws.Range("L132").Value = 62505948
$ws.Range("M132").Value = -1038.263
$ws.Range("N132").Value = -62511008
$ws.Range("H134").Value = 33334740
$ws.Range("I134").Value = 1267.0834
$ws.Range("J134").Value = 166668640
$ws.Range("K134").Value = 3801.2502
$ws.Range("L134").Value = 500005920
$ws.Range("M134").Value = -1266.2502
$ws.Range("N134").Value = -500010990
$ws.Range("H136").Value = 28572172
$ws.Range("I136").Value = 52632424
$ws.Range("J136").Value = 625.5
$ws.Range("K136").Value = 157897272
$ws.Range("L136").Value = 1876.5
$ws.Range("M136").Value = -157894722
$ws.Range("N136").Value = -6976.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 560.7143
$ws.Range("I34").Value = 292.85715
$ws.Range("J34").Value = 828.5714
$ws.Range("K34").Value = 878.5714499999999
$ws.Range("L34").Value = 2485.7142
$ws.Range("M34").Value = -794.5714499999999
$ws.Range("N34").Value = -2653.7142
$ws.Range("H122").Value = 14209809
$ws.Range("I122").Value = 125000350
$ws.Range("J122").Value = 5893.564
$ws.Range("K122").Value = 1125003150
$ws.Range("L122").Value = 53042.076
$ws.Range("M122").Value = -1125000700
$ws.Range("N122").Value = -57942.076
$ws.Range("H131").Value = 870.73
$ws.Range("J131").Value = 871.1414
$ws.Range("L131").Value = 2613.4242
$ws.Range("N131").Value = -12693.4242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 25184
$ws.Range("J116").Value = 25184
$ws.Range("L116").Value = 25184
$ws.Range("N116").Value = -34362
$ws.Range("H132").Value = 3660.0352
$ws.Range("I132").Value = 1634.766
$ws.Range("J132").Value = 13178.8
$ws.Range("K132").Value = 4904.298000000001
$ws.Range("L132").Value = 39536.39999999999
$ws.Range("M132").Value = -2374.298000000001
$ws.Range("N132").Value = -44596.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1056.25
$ws.Range("I32").Value = 1056.25
$ws.Range("K32").Value = 1056.25
$ws.Range("M32").Value = -739.25
$ws.Range("H132").Value = 22733998
$ws.Range("I132").Value = 45456284
$ws.Range("J132").Value = 11709.454
$ws.Range("K132").Value = 136368852
$ws.Range("L132").Value = 35128.362
$ws.Range("M132").Value = -136366322
$ws.Range("N132").Value = -40188.362
$ws.Range("H136").Value = 43419020
$ws.Range("I136").Value = 17638454
$ws.Range("J136").Value = 142858350
$ws.Range("K136").Value = 52915362
$ws.Range("L136").Value = 428575050
$ws.Range("M136").Value = -52912812
$ws.Range("N136").Value = -428580150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10441615
$ws.Range("I132").Value = 38109.855
$ws.Range("J132").Value = 25006522
$ws.Range("K132").Value = 114329.565
$ws.Range("L132").Value = 75019566
$ws.Range("M132").Value = -111799.565
$ws.Range("N132").Value = -75024626
$ws.Range("H136").Value = 20002616
$ws.Range("I136").Value = 45456330
$ws.Range("J136").Value = 3272.8572
$ws.Range("K136").Value = 136368990
$ws.Range("L136").Value = 9818.571599999999
$ws.Range("M136").Value = -136366440
$ws.Range("N136").Value = -14918.5716
